$d = $word.ActiveDocument

# The three section headings below were split across multiple runs
# ("1." / digit / " Title") and need to become a single run whose
# text is the concatenation, e.g. "1.3 Functions and Scope".
# Because the concatenated text is identical to the original text,
# this is purely a run-merge with no character changes.
#
# Technique: remove the leading "1." (first 2 characters) of the
# heading paragraph -- this merges the remaining runs (digit + rest)
# into a single clean run -- then re-insert "1." at the very start of
# the paragraph, which merges into that same run without disturbing
# its run properties or introducing spurious xml:space attributes.

$targets = @("1.3 Functions and Scope", "1.4 Arrays and Strings", "1.5  Pointers")

foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    foreach ($target in $targets) {
        if ($t -like "$target*") {
            $start = $para.Range.Start
            $prefix = $d.Range($start, $start + 2)
            if ($prefix.Text -eq "1.") {
                $prefix.Text = ""
                $insertPoint = $d.Range($start, $start)
                $insertPoint.InsertBefore("1.")
            }
        }
    }
}
